# Update "想去人数" (want-to-go count) values in column F
# for the sheets that contain this data: "展览" and "全部类型".
# Both sheets share identical rows for the events being updated.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 178
    4  = 12449
    5  = 1273
    6  = 144
    9  = 161
    14 = 128
    15 = 43
    17 = 3707
    22 = 53
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
